# Auto-generated Excel COM-interop script updating pricing/profit columns (H-N)
# across all 8 Leve sheets, per the scheduled-runner data refresh diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (137 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 455236.44
$ws.Range("I2").Value = 606399.9399999999
$ws.Range("J2").Value = 1746
$ws.Range("K2").Value = 606399.9399999999
$ws.Range("L2").Value = 1746
$ws.Range("M2").Value = -606286.9399999999
$ws.Range("N2").Value = -1972
$ws.Range("H9").Value = 402.72726
$ws.Range("I9").Value = 91.2
$ws.Range("K9").Value = 91.2
$ws.Range("M9").Value = 77.8
$ws.Range("H15").Value = 1946
$ws.Range("I15").Value = 1946
$ws.Range("K15").Value = 5838
$ws.Range("M15").Value = -5669
$ws.Range("H32").Value = 3299.6
$ws.Range("J32").Value = 3633
$ws.Range("L32").Value = 3633
$ws.Range("N32").Value = -4285
$ws.Range("H33").Value = 3180.4
$ws.Range("J33").Value = 683
$ws.Range("L33").Value = 683
$ws.Range("N33").Value = -1141
$ws.Range("H40").Value = 3768
$ws.Range("J40").Value = 4478.8335
$ws.Range("L40").Value = 4478.8335
$ws.Range("N40").Value = -4828.8335
$ws.Range("H46").Value = 3300
$ws.Range("I46").Value = 3500
$ws.Range("J46").Value = 2900
$ws.Range("K46").Value = 10500
$ws.Range("L46").Value = 8700
$ws.Range("M46").Value = -10381
$ws.Range("N46").Value = -8938
$ws.Range("H51").Value = 8879.799999999999
$ws.Range("I51").Value = 11159.8
$ws.Range("J51").Value = 6599.8
$ws.Range("K51").Value = 11159.8
$ws.Range("L51").Value = 6599.8
$ws.Range("M51").Value = -10675.8
$ws.Range("N51").Value = -7567.8
$ws.Range("H60").Value = 3300
$ws.Range("I60").Value = 3500
$ws.Range("J60").Value = 2900
$ws.Range("K60").Value = 10500
$ws.Range("L60").Value = 8700
$ws.Range("M60").Value = -10016
$ws.Range("N60").Value = -9668
$ws.Range("H62").Value = 13729.167
$ws.Range("I62").Value = 7930.909
$ws.Range("K62").Value = 7930.909
$ws.Range("M62").Value = -7306.909
$ws.Range("H65").Value = 13729.167
$ws.Range("I65").Value = 7930.909
$ws.Range("K65").Value = 39654.545
$ws.Range("M65").Value = -36534.545
$ws.Range("H76").Value = 5697.4
$ws.Range("I76").Value = 6748.5
$ws.Range("J76").Value = 4996.6665
$ws.Range("K76").Value = 6748.5
$ws.Range("L76").Value = 4996.6665
$ws.Range("M76").Value = -6433.5
$ws.Range("N76").Value = -5626.6665
$ws.Range("H79").Value = 5697.4
$ws.Range("I79").Value = 6748.5
$ws.Range("J79").Value = 4996.6665
$ws.Range("K79").Value = 6748.5
$ws.Range("L79").Value = 4996.6665
$ws.Range("M79").Value = -5656.5
$ws.Range("N79").Value = -7180.6665
$ws.Range("H86").Value = 167674.33
$ws.Range("I86").Value = 167674.33
$ws.Range("K86").Value = 167674.33
$ws.Range("M86").Value = -166551.33
$ws.Range("H89").Value = 167674.33
$ws.Range("I89").Value = 167674.33
$ws.Range("K89").Value = 838371.6499999999
$ws.Range("M89").Value = -832755.6499999999
$ws.Range("H98").Value = 2879.8462
$ws.Range("I98").Value = 2953.1667
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 2953.1667
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -1455.1667
$ws.Range("N98").Value = -4996
$ws.Range("H111").Value = 1413.8572
$ws.Range("I111").Value = 1224.25
$ws.Range("K111").Value = 3672.75
$ws.Range("M111").Value = -605.75
$ws.Range("H112").Value = 3117
$ws.Range("J112").Value = 3117
$ws.Range("L112").Value = 9351
$ws.Range("N112").Value = -11567
$ws.Range("H113").Value = 4522.304
$ws.Range("I113").Value = 4336.1177
$ws.Range("J113").Value = 5049.8335
$ws.Range("K113").Value = 4336.1177
$ws.Range("L113").Value = 5049.8335
$ws.Range("M113").Value = -1082.1177
$ws.Range("N113").Value = -11557.8335
$ws.Range("H122").Value = 2879.8462
$ws.Range("I122").Value = 2953.1667
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 8859.500100000001
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -6409.500100000001
$ws.Range("N122").Value = -10900
$ws.Range("H127").Value = 2616
$ws.Range("I127").Value = 2616
$ws.Range("K127").Value = 7848
$ws.Range("M127").Value = -2888
$ws.Range("H132").Value = 3643.862
$ws.Range("I132").Value = 3228.4583
$ws.Range("J132").Value = 5637.8
$ws.Range("K132").Value = 9685.374899999999
$ws.Range("L132").Value = 16913.4
$ws.Range("M132").Value = -7155.374899999999
$ws.Range("N132").Value = -21973.4
$ws.Range("H135").Value = 686.6667
$ws.Range("I135").Value = 568.9167
$ws.Range("K135").Value = 5120.2503
$ws.Range("M135").Value = -2585.2503
$ws.Range("H137").Value = 45657.434
$ws.Range("I137").Value = 78560.16
$ws.Range("K137").Value = 235680.48
$ws.Range("M137").Value = -233130.48
$ws.Range("H138").Value = 2190.275
$ws.Range("I138").Value = 1787.875
$ws.Range("J138").Value = 3799.875
$ws.Range("K138").Value = 5363.625
$ws.Range("L138").Value = 11399.625
$ws.Range("M138").Value = -223.625
$ws.Range("N138").Value = -21679.625
$ws.Range("H141").Value = 2114.6667
$ws.Range("I141").Value = 2114.6667
$ws.Range("K141").Value = 6344.000100000001
$ws.Range("M141").Value = -1164.000100000001

# --- Sheet: ARM (80 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1744.6
$ws.Range("I2").Value = 1304.9048
$ws.Range("J2").Value = 4053
$ws.Range("K2").Value = 1304.9048
$ws.Range("L2").Value = 4053
$ws.Range("M2").Value = -1191.9048
$ws.Range("N2").Value = -4279
$ws.Range("H4").Value = 699.5
$ws.Range("I4").Value = 699.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 699.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = ""
$ws.Range("N4").Value = -583.5
$ws.Range("H5").Value = 10297.3
$ws.Range("I5").Value = 14560.429
$ws.Range("J5").Value = 350
$ws.Range("K5").Value = 14560.429
$ws.Range("L5").Value = 350
$ws.Range("M5").Value = -14448.429
$ws.Range("N5").Value = -574
$ws.Range("H32").Value = 20391.508
$ws.Range("I32").Value = 20471.178
$ws.Range("K32").Value = 20471.178
$ws.Range("M32").Value = -20184.178
$ws.Range("H39").Value = 14000
$ws.Range("I39").Value = 14000
$ws.Range("K39").Value = 14000
$ws.Range("M39").Value = -13480
$ws.Range("H45").Value = 14145.667
$ws.Range("I45").Value = 11846.7
$ws.Range("K45").Value = 11846.7
$ws.Range("M45").Value = -11469.7
$ws.Range("H61").Value = 7562.921
$ws.Range("I61").Value = 7774.029
$ws.Range("J61").Value = 5100
$ws.Range("K61").Value = 7774.029
$ws.Range("L61").Value = 5100
$ws.Range("M61").Value = -7562.029
$ws.Range("N61").Value = -5524
$ws.Range("H74").Value = 29533.432
$ws.Range("I74").Value = 31061.086
$ws.Range("J74").Value = 2799.5
$ws.Range("K74").Value = 31061.086
$ws.Range("L74").Value = 2799.5
$ws.Range("M74").Value = -30187.086
$ws.Range("N74").Value = -4547.5
$ws.Range("H77").Value = 29533.432
$ws.Range("I77").Value = 31061.086
$ws.Range("J77").Value = 2799.5
$ws.Range("K77").Value = 155305.43
$ws.Range("L77").Value = 13997.5
$ws.Range("M77").Value = -150937.43
$ws.Range("N77").Value = -22733.5
$ws.Range("H102").Value = 2949.6667
$ws.Range("I102").Value = 2267.2222
$ws.Range("K102").Value = 2267.2222
$ws.Range("M102").Value = -645.2222000000002
$ws.Range("H116").Value = 1744.6
$ws.Range("I116").Value = 1304.9048
$ws.Range("J116").Value = 4053
$ws.Range("K116").Value = 1304.9048
$ws.Range("L116").Value = 4053
$ws.Range("M116").Value = 989.0952
$ws.Range("N116").Value = -8641
$ws.Range("H124").Value = 18498.75
$ws.Range("J124").Value = 18498.75
$ws.Range("L124").Value = 18498.75
$ws.Range("N124").Value = -28318.75
$ws.Range("H132").Value = 43802.76
$ws.Range("I132").Value = 47022.434
$ws.Range("K132").Value = 141067.302
$ws.Range("M132").Value = -138537.302
$ws.Range("H136").Value = 7562.921
$ws.Range("I136").Value = 7774.029
$ws.Range("J136").Value = 5100
$ws.Range("K136").Value = 23322.087
$ws.Range("L136").Value = 15300
$ws.Range("M136").Value = -20772.087
$ws.Range("N136").Value = -20400

# --- Sheet: BSM (69 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1744.6
$ws.Range("I3").Value = 1304.9048
$ws.Range("J3").Value = 4053
$ws.Range("K3").Value = 1304.9048
$ws.Range("L3").Value = 4053
$ws.Range("M3").Value = -1190.9048
$ws.Range("N3").Value = -4281
$ws.Range("H4").Value = 10297.3
$ws.Range("I4").Value = 14560.429
$ws.Range("J4").Value = 350
$ws.Range("K4").Value = 14560.429
$ws.Range("L4").Value = 350
$ws.Range("M4").Value = -14445.429
$ws.Range("N4").Value = -580
$ws.Range("H20").Value = 5125
$ws.Range("H22").Value = 1039.625
$ws.Range("I22").Value = 1068.7
$ws.Range("J22").Value = 991.1667
$ws.Range("K22").Value = 1068.7
$ws.Range("L22").Value = 991.1667
$ws.Range("M22").Value = -895.7
$ws.Range("N22").Value = -1337.1667
$ws.Range("H25").Value = 6013.1665
$ws.Range("I25").Value = 4812.6
$ws.Range("K25").Value = 4812.6
$ws.Range("M25").Value = -4577.6
$ws.Range("H80").Value = 850
$ws.Range("I80").Value = 850
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 850
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = 148
$ws.Range("H83").Value = 850
$ws.Range("I83").Value = 850
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 4250
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = 742
$ws.Range("H86").Value = 2272.879
$ws.Range("I86").Value = 1809.2963
$ws.Range("J86").Value = 4359
$ws.Range("K86").Value = 1809.2963
$ws.Range("L86").Value = 4359
$ws.Range("M86").Value = -686.2963
$ws.Range("N86").Value = -6605
$ws.Range("H89").Value = 2272.879
$ws.Range("I89").Value = 1809.2963
$ws.Range("J89").Value = 4359
$ws.Range("K89").Value = 9046.4815
$ws.Range("L89").Value = 21795
$ws.Range("M89").Value = -3430.4815
$ws.Range("N89").Value = -33027
$ws.Range("H94").Value = 3726.6155
$ws.Range("I94").Value = 2363.375
$ws.Range("J94").Value = 5907.8
$ws.Range("K94").Value = 2363.375
$ws.Range("L94").Value = 5907.8
$ws.Range("M94").Value = -1912.375
$ws.Range("N94").Value = -6809.8
$ws.Range("H105").Value = 2984.7896
$ws.Range("I105").Value = 2766.3635
$ws.Range("K105").Value = 2766.3635
$ws.Range("M105").Value = -1019.3635
$ws.Range("H134").Value = 2263.8635
$ws.Range("I134").Value = 2067.75
$ws.Range("K134").Value = 6203.25
$ws.Range("M134").Value = -3668.25

# --- Sheet: CRP (67 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1062.25
$ws.Range("I6").Value = 1042.5714
$ws.Range("K6").Value = 1042.5714
$ws.Range("M6").Value = -929.5714
$ws.Range("H7").Value = 181.42857
$ws.Range("I7").Value = 200.14285
$ws.Range("K7").Value = 200.14285
$ws.Range("M7").Value = -87.14285000000001
$ws.Range("H22").Value = 423.36365
$ws.Range("I22").Value = 284.83334
$ws.Range("J22").Value = 589.6
$ws.Range("K22").Value = 284.83334
$ws.Range("L22").Value = 589.6
$ws.Range("M22").Value = 65.16665999999998
$ws.Range("N22").Value = -1289.6
$ws.Range("H31").Value = 2747.7646
$ws.Range("I31").Value = 2528
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 2528
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = -2233
$ws.Range("N31").Value = -10590
$ws.Range("H34").Value = 2747.7646
$ws.Range("I34").Value = 2528
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 2528
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -2326
$ws.Range("N34").Value = -10404
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = ""
$ws.Range("N38").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = ""
$ws.Range("N46").Value = 0
$ws.Range("H58").Value = 49638
$ws.Range("I58").Value = 49638
$ws.Range("K58").Value = 49638
$ws.Range("M58").Value = -49435
$ws.Range("H94").Value = 1377
$ws.Range("J94").Value = 1252.6666
$ws.Range("L94").Value = 1252.6666
$ws.Range("N94").Value = -2154.6666
$ws.Range("H122").Value = 1886.091
$ws.Range("I122").Value = 1856.125
$ws.Range("J122").Value = 1966
$ws.Range("K122").Value = 5568.375
$ws.Range("L122").Value = 5898
$ws.Range("M122").Value = -3118.375
$ws.Range("N122").Value = -10798
$ws.Range("H132").Value = 1475.0526
$ws.Range("I132").Value = 1475.0526
$ws.Range("K132").Value = 4425.1578
$ws.Range("M132").Value = -1895.1578
$ws.Range("H134").Value = 92832.27
$ws.Range("I134").Value = 112795
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 338385
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -335850
$ws.Range("N134").Value = -14070
$ws.Range("H136").Value = 49638
$ws.Range("I136").Value = 49638
$ws.Range("K136").Value = 148914
$ws.Range("M136").Value = -146364

# --- Sheet: CUL (27 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1005495.7
$ws.Range("I4").Value = 712447.75
$ws.Range("K4").Value = 2137343.25
$ws.Range("M4").Value = -2137231.25
$ws.Range("H37").Value = 64238.777
$ws.Range("J37").Value = 64238.777
$ws.Range("L37").Value = 192716.331
$ws.Range("N37").Value = -192940.331
$ws.Range("H68").Value = 17524.5
$ws.Range("J68").Value = 19785.285
$ws.Range("L68").Value = 59355.855
$ws.Range("N68").Value = -60977.855
$ws.Range("H71").Value = 17524.5
$ws.Range("J71").Value = 19785.285
$ws.Range("L71").Value = 178067.565
$ws.Range("N71").Value = -186179.565
$ws.Range("H97").Value = 495.3
$ws.Range("I97").Value = 399
$ws.Range("J97").Value = 506
$ws.Range("K97").Value = 1197
$ws.Range("L97").Value = 1518
$ws.Range("M97").Value = -701
$ws.Range("N97").Value = -2510
$ws.Range("H113").Value = 834.8333
$ws.Range("J113").Value = 791.25
$ws.Range("L113").Value = 2373.75
$ws.Range("N113").Value = -6713.75

# --- Sheet: GSM (50 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 41492
$ws.Range("J15").Value = 41492
$ws.Range("L15").Value = 41492
$ws.Range("N15").Value = -42068
$ws.Range("H19").Value = 2000
$ws.Range("J19").Value = 2000
$ws.Range("L19").Value = 2000
$ws.Range("N19").Value = -2576
$ws.Range("H70").Value = 7419
$ws.Range("I70").Value = 7834.8335
$ws.Range("J70").Value = 6920
$ws.Range("K70").Value = 7834.8335
$ws.Range("L70").Value = 6920
$ws.Range("M70").Value = -7564.8335
$ws.Range("N70").Value = -7460
$ws.Range("H73").Value = 7419
$ws.Range("I73").Value = 7834.8335
$ws.Range("J73").Value = 6920
$ws.Range("K73").Value = 7834.8335
$ws.Range("L73").Value = 6920
$ws.Range("M73").Value = -6898.8335
$ws.Range("N73").Value = -8792
$ws.Range("H81").Value = 41492
$ws.Range("J81").Value = 41492
$ws.Range("L81").Value = 41492
$ws.Range("N81").Value = -43488
$ws.Range("H84").Value = 41492
$ws.Range("J84").Value = 41492
$ws.Range("L84").Value = 124476
$ws.Range("N84").Value = -134460
$ws.Range("H92").Value = 38373.375
$ws.Range("J92").Value = 38373.375
$ws.Range("L92").Value = 38373.375
$ws.Range("N92").Value = -42117.375
$ws.Range("H97").Value = 638.88
$ws.Range("I97").Value = 478
$ws.Range("K97").Value = 478
$ws.Range("M97").Value = 18
$ws.Range("H122").Value = 8645.704
$ws.Range("I122").Value = 7423.095
$ws.Range("K122").Value = 22269.285
$ws.Range("M122").Value = -19819.285
$ws.Range("H126").Value = 7832.643
$ws.Range("I126").Value = 7520
$ws.Range("K126").Value = 22560
$ws.Range("M126").Value = -20090
$ws.Range("H132").Value = 59954.55
$ws.Range("I132").Value = 72880.875
$ws.Range("K132").Value = 218642.625
$ws.Range("M132").Value = -216112.625

# --- Sheet: LTW (67 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12453.228
$ws.Range("I7").Value = 19429
$ws.Range("K7").Value = 19429
$ws.Range("M7").Value = -19317
$ws.Range("H22").Value = 2394.15
$ws.Range("I22").Value = 1181
$ws.Range("J22").Value = 3607.3
$ws.Range("K22").Value = 1181
$ws.Range("L22").Value = 3607.3
$ws.Range("M22").Value = -886
$ws.Range("N22").Value = -4197.3
$ws.Range("H27").Value = 2394.15
$ws.Range("I27").Value = 1181
$ws.Range("J27").Value = 3607.3
$ws.Range("K27").Value = 1181
$ws.Range("L27").Value = 3607.3
$ws.Range("M27").Value = -1074
$ws.Range("N27").Value = -3821.3
$ws.Range("H46").Value = 4212.421
$ws.Range("I46").Value = 1199.5
$ws.Range("J46").Value = 5603
$ws.Range("K46").Value = 1199.5
$ws.Range("L46").Value = 5603
$ws.Range("M46").Value = -1011.5
$ws.Range("N46").Value = -5979
$ws.Range("H61").Value = 2017.3658
$ws.Range("I61").Value = 1399.1082
$ws.Range("K61").Value = 1399.1082
$ws.Range("M61").Value = -1197.1082
$ws.Range("H68").Value = 3566.111
$ws.Range("I68").Value = 2882.5
$ws.Range("K68").Value = 2882.5
$ws.Range("M68").Value = -2133.5
$ws.Range("H71").Value = 3566.111
$ws.Range("I71").Value = 2882.5
$ws.Range("K71").Value = 14412.5
$ws.Range("M71").Value = -10668.5
$ws.Range("H100").Value = 3381.3572
$ws.Range("I100").Value = 2375
$ws.Range("J100").Value = 3783.9
$ws.Range("K100").Value = 2375
$ws.Range("L100").Value = 3783.9
$ws.Range("M100").Value = -1834
$ws.Range("N100").Value = -4865.9
$ws.Range("H113").Value = 2017.3658
$ws.Range("I113").Value = 1399.1082
$ws.Range("K113").Value = 1399.1082
$ws.Range("M113").Value = 770.8918000000001
$ws.Range("H126").Value = 12453.228
$ws.Range("I126").Value = 19429
$ws.Range("K126").Value = 58287
$ws.Range("M126").Value = -55817
$ws.Range("H127").Value = 157599.2
$ws.Range("J127").Value = 157599.2
$ws.Range("L127").Value = 157599.2
$ws.Range("N127").Value = -167519.2
$ws.Range("H132").Value = 20517.13
$ws.Range("I132").Value = 21896.984
$ws.Range("J132").Value = 6028.6665
$ws.Range("K132").Value = 65690.952
$ws.Range("L132").Value = 18085.9995
$ws.Range("M132").Value = -63160.952
$ws.Range("N132").Value = -23145.9995
$ws.Range("H136").Value = 2548.6897
$ws.Range("I136").Value = 2300.4285
$ws.Range("K136").Value = 6901.2855
$ws.Range("M136").Value = -4351.2855

# --- Sheet: WVR (48 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = ""
$ws.Range("N18").Value = 0
$ws.Range("H41").Value = 30499.334
$ws.Range("I41").Value = 33000
$ws.Range("J41").Value = 29249
$ws.Range("K41").Value = 33000
$ws.Range("L41").Value = 29249
$ws.Range("M41").Value = -32610
$ws.Range("N41").Value = -30029
$ws.Range("H81").Value = 3016.125
$ws.Range("I81").Value = 3116.5
$ws.Range("J81").Value = 2982.6667
$ws.Range("K81").Value = 6233
$ws.Range("L81").Value = 5965.3334
$ws.Range("M81").Value = -5172
$ws.Range("N81").Value = -8087.3334
$ws.Range("H84").Value = 3016.125
$ws.Range("I84").Value = 3116.5
$ws.Range("J84").Value = 2982.6667
$ws.Range("K84").Value = 31165
$ws.Range("L84").Value = 29826.667
$ws.Range("M84").Value = -25861
$ws.Range("N84").Value = -40434.667
$ws.Range("H107").Value = 2874.75
$ws.Range("I107").Value = 750
$ws.Range("K107").Value = 2250
$ws.Range("M107").Value = -330
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = ""
$ws.Range("N108").Value = 0
$ws.Range("H122").Value = 79908.30499999999
$ws.Range("I122").Value = 1981
$ws.Range("J122").Value = 339666
$ws.Range("K122").Value = 5943
$ws.Range("L122").Value = 1018998
$ws.Range("M122").Value = -3493
$ws.Range("N122").Value = -1023898
$ws.Range("H132").Value = 21533.24
$ws.Range("I132").Value = 22707.844
$ws.Range("K132").Value = 68123.53200000001
$ws.Range("M132").Value = -65593.53200000001
$ws.Range("H136").Value = 3795.2415
$ws.Range("J136").Value = 6252.25
$ws.Range("L136").Value = 18756.75
$ws.Range("N136").Value = -23856.75

Write-Output "Applied 545 cell updates across 8 sheets."
